$p = $ppt.ActivePresentation

# Slide 3: title "3.1 Code" -> "9.1 Code"
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "9.1 Code"

# Slide 4: title "2.1 Code" -> "9.1 Code"
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "9.1 Code"

# Slide 5: title "3.2 Verify" -> "9.2 Verify"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "9.2 Verify"

# Slide 6: title "3.2 Verify" -> "9.2 Verify"
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "9.2 Verify"
